$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.601.53"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Value = "1.630.09"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'212.69"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  +1.14%  "
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").Value = "'18.91"
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("E11").Value = "  +3.31%  "
$ws.Range("D12").Value = "1.857.64"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").Value = "1.637.07"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").Value = "'4.07"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").Value = "'0.524"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").Value = "26.583.68"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").Value = "'62.90"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("D20").Value = "'208.79"
$ws.Range("E20").Value = "  +3.92%  "
$ws.Range("D21").Value = "'4.28"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").Value = "'6.17"
$ws.Range("E23").Value = "  +2.27%  "
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").Value = "'146.62"
$ws.Range("E25").Value = "  +2.43%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("E28").Value = "  +4.43%  "
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").Value = "'0.0519"
$ws.Range("E30").Value = "  +3.76%  "
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").Value = "'3.24"
$ws.Range("E32").Value = "  +1.71%  "
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("D37").Value = "1.163.93"
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "'0.503"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("D42").Value = "'0.789"
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "1.767.55"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("D45").Value = "'92.44"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").Value = "'54.48"
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.56"
$ws.Range("E49").Value = "  +4.63%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.409"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("E51").Value = "  -0.11%  "
